# Insert a new data row at row 58 (pushing existing rows 58..178 down to 59..179)
# and populate it with a new weekly price observation for Mango / Vega Monumental
# Concepcion, matching the rest of the table's layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(58).Insert()

$ws.Cells.Item(58, 1).Value2 = 11
$ws.Cells.Item(58, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(58, 3).Value2 = "Bíobío"
$ws.Cells.Item(58, 4).Value2 = 45125
$ws.Cells.Item(58, 5).Value2 = 8
$ws.Cells.Item(58, 6).Value2 = "Fruta"
$ws.Cells.Item(58, 7).Value2 = 100108
$ws.Cells.Item(58, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(58, 9).Value2 = 100108002
$ws.Cells.Item(58, 10).Value2 = "Mango"
$ws.Cells.Item(58, 11).Value2 = "Sin especificar"
$ws.Cells.Item(58, 12).Value2 = "Primera"
$ws.Cells.Item(58, 13).Value2 = 100
$ws.Cells.Item(58, 14).Value2 = 7500
$ws.Cells.Item(58, 15).Value2 = 8000
$ws.Cells.Item(58, 16).Value2 = 7750
$ws.Cells.Item(58, 17).Value2 = "$/bandeja 4 kilos"
$ws.Cells.Item(58, 18).Value2 = "Brasil"
$ws.Cells.Item(58, 19).Value2 = 1938
$ws.Cells.Item(58, 20).Value2 = 4
